$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that still has the original (pre-edit) text.
# ------------------------------------------------------------------
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Um den Taster zu entprellen musste*") {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # New wording, still ending with the paragraph mark that Range.Text
    # reports (the existing trailing "werden:" sentence end becomes
    # "gestellt:" further down, so we just replace the whole sentence).
    $newText = "Um den Taster zu entprellen, haben wir in der Datei openhab.cfg den Debounce auf 10 ms gestellt:"
    $r.Text = $newText

    $start = $target.Range.Start

    # Boundaries (character offsets from the paragraph start) that line
    # up with the five runs the diff introduces:
    #   "Um den Taster zu entprellen,"            0  .. 28
    #   " haben wir"                              28 .. 38
    #   " in der Datei openhab.cfg de"            38 .. 66
    #   "n"                                       66 .. 67
    #   " Debounce auf 10 ms gestellt:"           67 .. 96
    $b1 = $start + 28
    $b2 = $start + 38
    $b3 = $start + 66
    $b4 = $start + 67
    $endOfText = $start + $newText.Length

    $run1 = $d.Range($start, $b1)
    $run2 = $d.Range($b1, $b2)
    $run3 = $d.Range($b2, $b3)
    $run4 = $d.Range($b3, $b4)
    $run5 = $d.Range($b4, $endOfText)

    # Forcing each sub-range's Italic property through a real transition
    # (on -> off) makes the engine keep the runs as separate <w:r>
    # elements once it re-lays the paragraph out, instead of silently
    # re-merging them back into a single run because their resulting
    # formatting is identical.
    foreach ($run in @($run1, $run2, $run3, $run4)) {
        $run.Font.Italic = $true
        $run.Font.Italic = $false
    }
}
